$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for rows 2-5 with repulled data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 2
